# Assignment 3 - heatMap.xlsx final commit
# Converts the "version-like" 4-digit numbers stored in column A (rows 3-76)
# and in the header row (row 77, columns B:BW) from numbers into their
# dotted textual form, e.g. 2025 -> "2.0.2.5", 4710 -> "4.7.1.0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values currently in A3:A76 (top to bottom).
$colA = @(2025,2026,2027,2030,2032,2040,2042,2060,2070,2080,2082,2084,2100,2102,2104,2200,2202,2300,2302,2304,2306,2400,2402,2500,2502,2504,3006,3008,3010,3012,3014,3016,3020,3022,3030,3034,3040,3042,3050,3052,3100,3110,4000,4002,4004,4100,4102,4104,4200,4202,4204,4208,4300,4302,4304,4306,4310,4312,4314,4400,4402,4404,4406,4410,4500,4502,4504,4510,4600,4602,4604,4700,4702,4710)

# Values currently in B77:BW77 (left to right).
$row77 = @(2008,2010,2025,2026,2027,2030,2032,2040,2042,2060,2070,2080,2082,2084,2100,2102,2104,2200,2202,2300,2302,2304,2306,2400,2402,2500,2502,2504,3006,3008,3010,3012,3014,3016,3020,3022,3030,3034,3040,3042,3050,3052,3100,3110,4000,4002,4004,4100,4102,4104,4200,4202,4204,4208,4300,4302,4304,4306,4310,4312,4314,4400,4402,4404,4406,4410,4500,4502,4504,4510,4600,4602,4604,4700)

# Column A first (this is the order the workbook's shared-string table was
# built in: A3..A76 top-to-bottom, then B77..BW77 left-to-right).
for ($i = 0; $i -lt $colA.Length; $i++) {
    $n = $colA[$i].ToString().PadLeft(4, '0')
    $text = "{0}.{1}.{2}.{3}" -f $n.Substring(0,1), $n.Substring(1,1), $n.Substring(2,1), $n.Substring(3,1)
    $row = $i + 3
    $ws.Cells.Item($row, 1).Value = $text
}

for ($i = 0; $i -lt $row77.Length; $i++) {
    $n = $row77[$i].ToString().PadLeft(4, '0')
    $text = "{0}.{1}.{2}.{3}" -f $n.Substring(0,1), $n.Substring(1,1), $n.Substring(2,1), $n.Substring(3,1)
    $col = $i + 2
    $ws.Cells.Item(77, $col).Value = $text
}
